$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25, columns B-F and I-M (case with 380 kV)
$data = @(
    @(1.02, 1.041804464500799, 1.051057722846995, 1.0400258124473, 1.059670161041093, 1.046515349916376, 1.046883698422676, 1.053810279920723, 1.0428093309011, 1.062399032429219),
    @(1.02, 1.042699779461084, 1.051785684973363, 1.040786244702253, 1.060529093033269, 1.046786279295741, 1.047425461809877, 1.05435103014468, 1.043380233574392, 1.06307212643042),
    @(1.02, 1.043279636703859, 1.05225715070812, 1.041279114036584, 1.061085641440233, 1.046960567263924, 1.047775896780013, 1.054700686201351, 1.043749801297007, 1.063507779662342),
    @(1.02, 1.043523533902915, 1.052455455200363, 1.041486510733213, 1.061319795130717, 1.047033592648053, 1.04792318925728, 1.0548476216312, 1.043905203574987, 1.063690954743024),
    @(1.02, 1.0435644926198, 1.052488757232649, 1.04152134490608, 1.061359121115607, 1.047045839526733, 1.047947918497557, 1.054872289184393, 1.043931298377179, 1.063721712155978),
    @(1.02, 1.043282895179825, 1.052259800070339, 1.041281884519431, 1.061088769505224, 1.046961543997376, 1.047777865030116, 1.054702649796931, 1.043751877649902, 1.063510227154519),
    @(1.02, 1.042106930131979, 1.051303651989023, 1.040282633586254, 1.059960282349027, 1.046607122823236, 1.047066814507837, 1.053993079098522, 1.043002237195611, 1.062626482750753),
    @(1.02, 1.040038836832393, 1.049622126129068, 1.038528167532488, 1.057977653747417, 1.045974796983972, 1.045812967223023, 1.052740897854564, 1.041682520098875, 1.061070165836893),
    @(1.02, 1.038662948291843, 1.048503441270084, 1.037362882314298, 1.056659974472932, 1.045548052396475, 1.044976541133953, 1.051904953384449, 1.040803618527564, 1.060033345802536),
    @(1.02, 1.038067863425091, 1.048019611322685, 1.036859354773078, 1.056090391787378, 1.045362044470579, 1.044614248293596, 1.051542720730963, 1.040423275478066, 1.05958458061528),
    @(1.02, 1.037846926454083, 1.047839982251595, 1.036672481379919, 1.055878972495927, 1.045292769557615, 1.044479660468146, 1.051408133040308, 1.040282034474446, 1.059417918388518),
    @(1.02, 1.037894313455375, 1.047878509338685, 1.036712559129195, 1.055924315878616, 1.045307637544308, 1.044508530725691, 1.051437004273414, 1.040312329530453, 1.05945366666068),
    @(1.02, 1.038049598576195, 1.048004761345471, 1.03684390450409, 1.05607291275195, 1.045356321919549, 1.044603123544608, 1.051531596446035, 1.040411599724693, 1.059570803661023),
    @(1.02, 1.038145288689747, 1.048082560937927, 1.036924851862195, 1.05616448799438, 1.045386293703876, 1.044661403197912, 1.051589872758619, 1.040472768086929, 1.059642979536586),
    @(1.02, 1.038702456634549, 1.048535563571906, 1.03739632200994, 1.056697796614182, 1.045560371387419, 1.04500058299274, 1.051928988128898, 1.040828865532294, 1.060063132880839),
    @(1.02, 1.039052137387636, 1.048819873336705, 1.037692344665592, 1.057032590881072, 1.045669238332236, 1.04521331156393, 1.052141636616666, 1.04105229771943, 1.060326734412798),
    @(1.02, 1.039256166159455, 1.048985761108059, 1.037865110743881, 1.057227965265747, 1.045732620365364, 1.045337381386543, 1.052265645409113, 1.041182643767281, 1.060480506456965),
    @(1.02, 1.0393257458375, 1.04904233378038, 1.037924036606467, 1.057294598868912, 1.045754211928012, 1.045379684063928, 1.052307924845998, 1.041227092067468, 1.060532941715256),
    @(1.02, 1.039014613126025, 1.048789363920767, 1.037660573770255, 1.056996660832534, 1.045657570161648, 1.045190488947106, 1.052118824072614, 1.04102832329319, 1.06029845061913),
    @(1.02, 1.03800386809598, 1.047967580864273, 1.036805222158511, 1.056029150547686, 1.045341990638374, 1.044575268753482, 1.0515037424581, 1.040382366141418, 1.059536308918259),
    @(1.02, 1.037368975190797, 1.047451396879706, 1.036268350269637, 1.055421701783376, 1.045142513176553, 1.044188362610633, 1.051116794913298, 1.039976432309079, 1.059057289495434),
    @(1.02, 1.037705486380272, 1.047724987509291, 1.036552868320496, 1.055743639471677, 1.045248360199719, 1.044393477299316, 1.051321943710737, 1.040191605670056, 1.059311210240929),
    @(1.02, 1.039031568521126, 1.048803149643752, 1.037674929358194, 1.057012895782918, 1.04566284287102, 1.045200801539846, 1.052129132158154, 1.041039156237118, 1.060311230793398),
    @(1.02, 1.040572993677316, 1.050056436651793, 1.038980977741183, 1.05848950014274, 1.046139187070061, 1.046137214785242, 1.053064825876424, 1.042023543558614, 1.061472389461825)
)

$colIndices = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $ws.Cells.Item($rowNum, $colIndices[$j]).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($data.Count) rows"